$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly rental-price data rows for Sep-Dec 2025 (columns A:U)
$data = @(
    @(2025, 9, 12.4, 10, 9.9, 19.2, 15, 11.5, 8.9, 8, 19.1, 12.4, 7.3, 9.3, 20.9, 8.9, 10.7, 14.9, 9, 13.9, 10.2),
    @(2025, 10, 12.4, 10.1, 9.9, 19.1, 15.1, 11.1, 9, 8, 19.1, 12.3, 7.2, 9.1, 21.1, 8.8, 10.6, 14.9, 9, 13.9, 10.3),
    @(2025, 11, 12.7, 10.2, 10, 19.1, 15.1, 10.8, 9.1, 8.1, 18.7, 12.4, 7.2, 9.1, 21, 8.8, 10.7, 14.8, 8.9, 13.9, 10.3),
    @(2025, 12, 13, 10.4, 10.2, 19.1, 15.3, 10.9, 9.3, 8.2, 18.6, 12.5, 7.3, 9.4, 20.8, 9, 11.1, 14.9, 9.1, 13.9, 10.4)
)

$startRow = 238
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

# Scroll the view down to the newly added rows and move the selection
# to match where the author left off editing.
$win = $excel.ActiveWindow
$win.ScrollRow = 222
$win.ScrollColumn = 1
$ws.Range("V240").Select() | Out-Null
